$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$left = New-Object "object[,]" 24,5
$left[0,0] = 1.02
$left[0,1] = 1.051448455994865
$left[0,2] = 1.057778894763602
$left[0,3] = 1.062060961295877
$left[0,4] = 1.070434443336732
$left[1,0] = 1.02
$left[1,1] = 1.052415881977283
$left[1,2] = 1.05854435060914
$left[1,3] = 1.062993405530945
$left[1,4] = 1.071352091332666
$left[2,0] = 1.02
$left[2,1] = 1.053042139583332
$left[2,2] = 1.059039910001685
$left[2,3] = 1.063597842955952
$left[2,4] = 1.071946631557209
$left[3,0] = 1.02
$left[3,1] = 1.053305481739483
$left[3,2] = 1.059248303819954
$left[3,3] = 1.063852206578993
$left[3,4] = 1.072196757093486
$left[4,0] = 1.02
$left[4,1] = 1.053349701745097
$left[4,2] = 1.059283297581186
$left[4,3] = 1.063894930449665
$left[4,4] = 1.072238764838593
$left[5,0] = 1.02
$left[5,1] = 1.053045658125721
$left[5,2] = 1.059042694333256
$left[5,3] = 1.063601240762487
$left[5,4] = 1.071949973037929
$left[6,0] = 1.02
$left[6,1] = 1.051775346114019
$left[6,2] = 1.058037529962028
$left[6,3] = 1.062375860031313
$left[6,4] = 1.070744408522773
$left[7,0] = 1.02
$left[7,1] = 1.049538983046836
$left[7,2] = 1.056268332115336
$left[7,3] = 1.060224939509735
$left[7,4] = 1.068625931034483
$left[8,0] = 1.02
$left[8,1] = 1.048049523141602
$left[8,2] = 1.055090301239965
$left[8,3] = 1.058796687757337
$left[8,4] = 1.067217641244279
$left[9,0] = 1.02
$left[9,1] = 1.047404922974612
$left[9,2] = 1.054580555411641
$left[9,3] = 1.058179605707091
$left[9,4] = 1.066608807097382
$left[10,0] = 1.02
$left[10,1] = 1.047165542360352
$left[10,2] = 1.054391266715083
$left[10,3] = 1.057950599595523
$left[10,4] = 1.06638280511494
$left[11,0] = 1.02
$left[11,1] = 1.047216887924705
$left[11,2] = 1.054431867339532
$left[11,3] = 1.057999712858135
$left[11,4] = 1.066431276675409
$left[12,0] = 1.02
$left[12,1] = 1.047385134621874
$left[12,2] = 1.054564907649022
$left[12,3] = 1.058160671788261
$left[12,4] = 1.066590122715074
$left[13,0] = 1.02
$left[13,1] = 1.047488803952453
$left[13,2] = 1.054646885310678
$left[13,3] = 1.058259871183967
$left[13,4] = 1.066688012402123
$left[14,0] = 1.02
$left[14,1] = 1.04809231047154
$left[14,2] = 1.055124138847032
$left[14,3] = 1.058837670229525
$left[14,4] = 1.067258067987747
$left[15,0] = 1.02
$left[15,1] = 1.048470967194731
$left[15,2] = 1.055423601693548
$left[15,3] = 1.059200473404473
$left[15,4] = 1.067615907903098
$left[16,0] = 1.02
$left[16,1] = 1.048691864576885
$left[16,2] = 1.055598306934092
$left[16,3] = 1.059412221531301
$left[16,4] = 1.067824722863115
$left[17,0] = 1.02
$left[17,1] = 1.04876719050839
$left[17,2] = 1.055657882576265
$left[17,3] = 1.059484444381543
$left[17,4] = 1.067895939104592
$left[18,0] = 1.02
$left[18,1] = 1.048430337460616
$left[18,2] = 1.055391468677715
$left[18,3] = 1.059161534480407
$left[18,4] = 1.067577505439466
$left[19,0] = 1.02
$left[19,1] = 1.047335588706162
$left[19,2] = 1.054525729105188
$left[19,3] = 1.058113267708445
$left[19,4] = 1.066543342471999
$left[20,0] = 1.02
$left[20,1] = 1.046647581293612
$left[20,2] = 1.053981714945501
$left[20,3] = 1.057455371161004
$left[20,4] = 1.065893969217595
$left[21,0] = 1.02
$left[21,1] = 1.047012277956251
$left[21,2] = 1.054270077322625
$left[21,3] = 1.057804021314971
$left[21,4] = 1.066238133668197
$left[22,0] = 1.02
$left[22,1] = 1.048448696185613
$left[22,2] = 1.055405988100611
$left[22,3] = 1.059179128898901
$left[22,4] = 1.067594857572542
$left[23,0] = 1.02
$left[23,1] = 1.050116884400561
$left[23,2] = 1.05672546442932
$left[23,3] = 1.060780005804429
$left[23,4] = 1.069172903635301

$right = New-Object "object[,]" 24,6
$right[0,0] = 1.051412423241218
$right[0,1] = 1.056476274376347
$right[0,2] = 1.060512931665655
$right[0,3] = 1.064783331318106
$right[0,4] = 1.073134290773742
$right[0,5] = 1.022606431467129
$right[1,0] = 1.051716714334461
$right[1,1] = 1.057093700422855
$right[1,2] = 1.061092413891692
$right[1,3] = 1.065530221841768
$right[1,4] = 1.07386804928142
$right[1,5] = 1.022815308762946
$right[2,0] = 1.051912489705265
$right[2,1] = 1.057492848527447
$right[2,2] = 1.061466972259166
$right[2,3] = 1.066013962190819
$right[2,4] = 1.074342965428181
$right[2,5] = 1.022950259743203
$right[3,0] = 1.05199452470393
$right[3,1] = 1.057660561522775
$right[3,2] = 1.061624338691483
$right[3,3] = 1.066217433822144
$right[3,4] = 1.074542649130929
$right[3,5] = 1.023006943323905
$right[4,0] = 1.05200828295367
$right[4,1] = 1.057688716037972
$right[4,2] = 1.061650755441839
$right[4,3] = 1.066251603858657
$right[4,4] = 1.07457617857073
$right[4,5] = 1.023016457815348
$right[5,0] = 1.051913586918438
$right[5,1] = 1.057495089866897
$right[5,2] = 1.061469075382909
$right[5,3] = 1.066016680568886
$right[5,4] = 1.074345633499131
$right[5,5] = 1.022951017348148
$right[6,0] = 1.051515491701684
$right[6,1] = 1.05668501207587
$right[6,2] = 1.060708853818954
$right[6,3] = 1.065035652180317
$right[6,4] = 1.073382241315964
$right[6,5] = 1.02267706509313
$right[7,0] = 1.050805432196789
$right[7,1] = 1.055254769939895
$right[7,2] = 1.059366179300555
$right[7,3] = 1.063310466010015
$right[7,4] = 1.071685627341684
$right[7,5] = 1.022192759428425
$right[8,0] = 1.050326331103516
$right[8,1] = 1.054299451999676
$right[8,2] = 1.05846905089421
$right[8,3] = 1.062162766034779
$right[8,4] = 1.070555289335952
$right[8,5] = 1.021868856154031
$right[9,0] = 1.050117522608667
$right[9,1] = 1.053885365921842
$right[9,2] = 1.058080117408099
$right[9,3] = 1.06166638737284
$right[9,4] = 1.070066028213932
$right[9,5] = 1.021728361725186
$right[10,0] = 1.050039758785563
$right[10,1] = 1.05373149247045
$right[10,2] = 1.057935580219597
$right[10,3] = 1.061482098853895
$right[10,4] = 1.06988432327932
$right[10,5] = 1.021676139780732
$right[11,0] = 1.050056448570722
$right[11,1] = 1.053764501707497
$right[11,2] = 1.057966587091083
$right[11,3] = 1.061521625340145
$right[11,4] = 1.069923298302271
$right[11,5] = 1.021687343191559
$right[12,0] = 1.050111098766318
$right[12,1] = 1.053872647984417
$right[12,2] = 1.058068171341934
$right[12,3] = 1.061651152212655
$right[12,4] = 1.07005100784169
$right[12,5] = 1.021724045776654
$right[13,0] = 1.050144743662992
$right[13,1] = 1.053939272089504
$right[13,2] = 1.058130751502735
$right[13,3] = 1.06173096976718
$right[13,4] = 1.070129697700618
$right[13,5] = 1.021746654682379
$right[14,0] = 1.050340160535501
$right[14,1] = 1.054326924591932
$right[14,2] = 1.05849485327283
$right[14,3] = 1.06219572142266
$right[14,4] = 1.07058776391827
$right[14,5] = 1.021878175232047
$right[15,0] = 1.050462378044031
$right[15,1] = 1.054569974926464
$right[15,2] = 1.058723119154949
$right[15,3] = 1.062487404668706
$right[15,4] = 1.070875146360196
$right[15,5] = 1.021960609991268
$right[16,0] = 1.050533534722867
$right[16,1] = 1.054711700816426
$right[16,2] = 1.058856217288062
$right[16,3] = 1.062657594753341
$right[16,4] = 1.071042789237721
$right[16,5] = 1.022008669403531
$right[17,0] = 1.050557775124229
$right[17,1] = 1.054760018662329
$right[17,2] = 1.05890159257153
$right[17,3] = 1.062715634672303
$right[17,4] = 1.071099954075498
$right[17,5] = 1.022025052444211
$right[18,0] = 1.050449278780323
$right[18,1] = 1.054543902188036
$right[18,2] = 1.058698633086057
$right[18,3] = 1.062456103998241
$right[18,4] = 1.070844311120366
$right[18,5] = 1.021951767938537
$right[19,0] = 1.050095011247356
$right[19,1] = 1.053840803346912
$right[19,2] = 1.058038259237322
$right[19,3] = 1.061613007322535
$right[19,4] = 1.070013399772686
$right[19,5] = 1.021713238772892
$right[20,0] = 1.049871094626449
$right[20,1] = 1.053398370486458
$right[20,2] = 1.05762265152641
$right[20,3] = 1.061083432203097
$right[20,4] = 1.069491137807406
$right[20,5] = 1.021563057559959
$right[21,0] = 1.049989908245658
$right[21,1] = 1.053632947035103
$right[21,2] = 1.057843011114013
$right[21,3] = 1.061364121002436
$right[21,4] = 1.069767982790075
$right[21,5] = 1.021642691113115
$right[22,0] = 1.050455198177683
$right[22,1] = 1.054555683463762
$right[22,2] = 1.058709697427219
$right[22,3] = 1.062470247250399
$right[22,4] = 1.070858244183772
$right[22,5] = 1.021955763353854
$right[23,0] = 1.050990011451322
$right[23,1] = 1.055624846609723
$right[23,2] = 1.059713651153187
$right[23,3] = 1.063756045365187
$right[23,4] = 1.072124117156851
$right[23,5] = 1.022318147603484

$ws.Range("B2:F25").Value = $left
$ws.Range("I2:N25").Value = $right

Write-Output "Updated vm_pu values for 380 kV case"
